$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row at position 23 (shifts existing rows 23-40 down to 24-41)
$ws1.Rows.Item(23).Insert()

# Populate the new row with the new client "JUNCO SANCHEZ ARTURO ENRIQUE"
$ws1.Range("A23").Value = "OFICINA-CATAECSA"
$ws1.Range("B23").Value = "JUNCO SANCHEZ ARTURO ENRIQUE"
for ($col = 3; $col -le 18; $col++) {
  $ws1.Cells.Item(23, $col).Value = 0
}

# Update the summary/count row (now row 41) from "de 38" to "de 39"
$ws1.Range("C41").Value = "0 de 39"
$ws1.Range("D41").Value = "1 de 39"
$ws1.Range("E41").Value = "1 de 39"
$ws1.Range("F41").Value = "0 de 39"
$ws1.Range("G41").Value = "0 de 39"
$ws1.Range("H41").Value = "0 de 39"
$ws1.Range("I41").Value = "0 de 39"
$ws1.Range("J41").Value = "0 de 39"
$ws1.Range("K41").Value = "0 de 39"
$ws1.Range("L41").Value = "3 de 39"
$ws1.Range("M41").Value = "3 de 39"
$ws1.Range("N41").Value = "0 de 39"
$ws1.Range("O41").Value = "0 de 39"
$ws1.Range("P41").Value = "0 de 39"
$ws1.Range("Q41").Value = "0 de 39"
$ws1.Range("R41").Value = "0 de 39"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a new row at position 23 (shifts existing rows 23-40 down to 24-41)
$ws2.Rows.Item(23).Insert()

# Populate the new row with the new client "JUNCO SANCHEZ ARTURO ENRIQUE"
$ws2.Range("A23").Value = "OFICINA-CATAECSA"
$ws2.Range("B23").Value = "JUNCO SANCHEZ ARTURO ENRIQUE"
for ($col = 3; $col -le 7; $col++) {
  $ws2.Cells.Item(23, $col).Value = 0
}
